$d = $word.ActiveDocument

# The document's Word package XML (pkg:package) containing every part,
# including headers/footers, as a single string we can splice.
$xml = $d.XML

function Replace-InPart($xml, $partName, $oldText, $newText) {
    $marker = '<pkg:part pkg:name="/word/' + $partName + '"'
    $partStart = $xml.IndexOf($marker)
    if ($partStart -lt 0) {
        return $xml
    }
    $partEndMarker = "</pkg:part>"
    $partEndIdx = $xml.IndexOf($partEndMarker, $partStart)
    $partEnd = $partEndIdx + $partEndMarker.Length

    $before = $xml.Substring(0, $partStart)
    $part = $xml.Substring($partStart, $partEnd - $partStart)
    $after = $xml.Substring($partEnd)

    $part = $part.Replace($oldText, $newText)

    return $before + $part + $after
}

# Pearson logo pictures in the footers: image1.png -> image2.png
$xml = Replace-InPart $xml "footer1.xml" ' name="image1.png"' ' name="image2.png"'
$xml = Replace-InPart $xml "footer2.xml" ' name="image1.png"' ' name="image2.png"'

# BTEC logo pictures in the headers: image2.jpg -> image1.jpg
$xml = Replace-InPart $xml "header1.xml" ' name="image2.jpg"' ' name="image1.jpg"'
$xml = Replace-InPart $xml "header2.xml" ' name="image2.jpg"' ' name="image1.jpg"'

$d.XML = $xml
